# Update EUR->ARS rate: 2025-10-11T15:17:42Z
# Append the next data row (71) to the rate-history sheet.
#
# Columns A (date) and B (time) are plain text that happen to look like a
# date/time, so a leading apostrophe forces Excel to store them as literal
# text instead of auto-converting to a date/time serial number. The
# apostrophe itself is not stored in the cell's value. Re-applying the
# "Normal" style afterwards clears the quote-prefix formatting flag that
# Excel would otherwise stamp on the cell, keeping it identical in style to
# all the other rows in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A71").Value = "'2025-10-11"
$ws.Range("B71").Value = "'15:17:42"
$ws.Range("C71").Value = "1.00 EUR = 1,756.2048"

$ws.Range("A71:C71").Style = "Normal"
